$d = $word.ActiveDocument
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.ParagraphFormat.Alignment = 1
    $p.Range.Font.Name = "Times New Roman"
    $p.Range.Font.Size = 12
}
Write-Output "done"
